# Updates cryptos list D (Price) / E (Volume 1h) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.861.62"
$ws.Range("E2").Value = "  -0.90%  "

$ws.Range("D3").Value = "2.367.59"
$ws.Range("E3").Value = "  -1.14%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.636"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.36%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.34%  "

$ws.Range("E11").Value = "  -1.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("E13").Value = "  -4.87%  "

$ws.Range("E14").Value = "  -0.05%  "

$ws.Range("E15").Value = "  -5.26%  "

$ws.Range("D16").Value = "2.724.94"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "2.395.18"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "42.808.33"
$ws.Range("E18").Value = "  -1.00%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("E20").Value = "  -1.75%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("E22").Value = "  -1.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -6.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.03%  "

$ws.Range("E25").Value = "  -2.11%  "

$ws.Range("E26").Value = "  +0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.69%  "

$ws.Range("E29").Value = "  +2.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0897"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.65%  "

$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.02%  "

$ws.Range("E35").Value = "  +12.15%  "

$ws.Range("E36").Value = "  -2.43%  "

$ws.Range("E37").Value = "  -3.68%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0366"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.94"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.10%  "

$ws.Range("E40").Value = "  -4.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.242"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.68%  "

$ws.Range("E42").Value = "  -4.41%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.31%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.27%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.94%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "113.14"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.78%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "77.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.88%  "
